$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update existing header labels ---
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" worksheet after "Monthly Trend" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# --- Header row ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- Copy header formatting (bold, border, centered) from Weekly Quantity sheet ---
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# --- Copy date-number formatting for column A from Weekly Quantity sheet ---
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Forecast data rows ---
$newSheet.Range("A2").Value = 45200.99999999999
$newSheet.Range("B2").Value = 0
$newSheet.Range("C2").Value = -105.0817688391891
$newSheet.Range("D2").Value = 55.57197027014781
$newSheet.Range("A3").Value = 45207.99999999999
$newSheet.Range("B3").Value = 0
$newSheet.Range("C3").Value = -87.29919717297663
$newSheet.Range("D3").Value = 56.82692180663594
$newSheet.Range("A4").Value = 45214.99999999999
$newSheet.Range("B4").Value = 0
$newSheet.Range("C4").Value = -82.84333970974471
$newSheet.Range("D4").Value = 70.37506351786
$newSheet.Range("A5").Value = 45221.99999999999
$newSheet.Range("B5").Value = 3
$newSheet.Range("C5").Value = -72.92769269847756
$newSheet.Range("D5").Value = 82.59687695038444
$newSheet.Range("A6").Value = 45228.99999999999
$newSheet.Range("B6").Value = 12
$newSheet.Range("C6").Value = -68.64948363056973
$newSheet.Range("D6").Value = 89.15982893508608
$newSheet.Range("A7").Value = 45235.99999999999
$newSheet.Range("B7").Value = 21
$newSheet.Range("C7").Value = -59.79496614587347
$newSheet.Range("D7").Value = 95.26710515995396
$newSheet.Range("A8").Value = 45242.99999999999
$newSheet.Range("B8").Value = 30
$newSheet.Range("C8").Value = -47.58147224157027
$newSheet.Range("D8").Value = 103.218490489506
$newSheet.Range("A9").Value = 45249.99999999999
$newSheet.Range("B9").Value = 39
$newSheet.Range("C9").Value = -38.47385806126292
$newSheet.Range("D9").Value = 113.1124256320823
$newSheet.Range("A10").Value = 45256.99999999999
$newSheet.Range("B10").Value = 48
$newSheet.Range("C10").Value = -30.39301938666001
$newSheet.Range("D10").Value = 118.8109489159688
$newSheet.Range("A11").Value = 45263.99999999999
$newSheet.Range("B11").Value = 57
$newSheet.Range("C11").Value = -24.28162228495999
$newSheet.Range("D11").Value = 130.7383221802064
$newSheet.Range("A12").Value = 45270.99999999999
$newSheet.Range("B12").Value = 66
$newSheet.Range("C12").Value = -13.44031377914871
$newSheet.Range("D12").Value = 139.3332021469703
$newSheet.Range("A13").Value = 45277.99999999999
$newSheet.Range("B13").Value = 75
$newSheet.Range("C13").Value = -4.449269261191917
$newSheet.Range("D13").Value = 149.1929894262678
$newSheet.Range("A14").Value = 45298.99999999999
$newSheet.Range("B14").Value = 102
$newSheet.Range("C14").Value = 16.49830357590232
$newSheet.Range("D14").Value = 178.5876941096271
$newSheet.Range("A15").Value = 45305.99999999999
$newSheet.Range("B15").Value = 111
$newSheet.Range("C15").Value = 31.0779089696858
$newSheet.Range("D15").Value = 186.7075995213061
$newSheet.Range("A16").Value = 45312.99999999999
$newSheet.Range("B16").Value = 120
$newSheet.Range("C16").Value = 38.92660745673115
$newSheet.Range("D16").Value = 198.5429615197041
$newSheet.Range("A17").Value = 45319.99999999999
$newSheet.Range("B17").Value = 129
$newSheet.Range("C17").Value = 53.37416387443196
$newSheet.Range("D17").Value = 203.063791138876
$newSheet.Range("A18").Value = 45326.99999999999
$newSheet.Range("B18").Value = 138
$newSheet.Range("C18").Value = 55.37170684517007
$newSheet.Range("D18").Value = 216.6512783176417
$newSheet.Range("A19").Value = 45333.99999999999
$newSheet.Range("B19").Value = 147
$newSheet.Range("C19").Value = 72.40703070496127
$newSheet.Range("D19").Value = 225.1530074870252
$newSheet.Range("A20").Value = 45340.99999999999
$newSheet.Range("B20").Value = 156
$newSheet.Range("C20").Value = 75.27798815596988
$newSheet.Range("D20").Value = 232.8367281213133
$newSheet.Range("A21").Value = 45347.99999999999
$newSheet.Range("B21").Value = 165
$newSheet.Range("C21").Value = 87.63840388313643
$newSheet.Range("D21").Value = 241.4402714921206
$newSheet.Range("A22").Value = 45354.99999999999
$newSheet.Range("B22").Value = 174
$newSheet.Range("C22").Value = 91.90955841492384
$newSheet.Range("D22").Value = 248.2851659814723
$newSheet.Range("A23").Value = 45361.99999999999
$newSheet.Range("B23").Value = 183
$newSheet.Range("C23").Value = 107.0325003839219
$newSheet.Range("D23").Value = 256.7033729401662
$newSheet.Range("A24").Value = 45368.99999999999
$newSheet.Range("B24").Value = 191
$newSheet.Range("C24").Value = 114.2342679110252
$newSheet.Range("D24").Value = 266.7741524753727
$newSheet.Range("A25").Value = 45375.99999999999
$newSheet.Range("B25").Value = 200
$newSheet.Range("C25").Value = 125.3723880172993
$newSheet.Range("D25").Value = 271.9898787109265
$newSheet.Range("A26").Value = 45382.99999999999
$newSheet.Range("B26").Value = 209
$newSheet.Range("C26").Value = 127.429528623238
$newSheet.Range("D26").Value = 287.6357542479714
$newSheet.Range("A27").Value = 45389.99999999999
$newSheet.Range("B27").Value = 218
$newSheet.Range("C27").Value = 142.1514945165898
$newSheet.Range("D27").Value = 294.7661806015753
